# Update gh-pages to output generated at 456a3b4
#
# The scraped feed dropped "2024-10-05 nanning hua hai" event; every
# later row shifted up one slot, the stale duplicate last row was
# removed, and two 'want to go' counts (F) ticked up. Applies to the
# "展览" (sheet 1) and "全部类型" (sheet 4) tabs only.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("B2").Value = "'2024-10-19"
$ws1.Range("C2").Value = "南宁·10.19剑网3同人only——寒光烈火·阵营PK战"
$ws1.Range("D2").Value = "大学东路158号 维也纳酒店动物园店"
$ws1.Range("E2").Value = "2024.10.19 10:00-10.19 17:30"
$ws1.Range("F2").Value = 7
$ws1.Range("G2").Value = 78
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=92730"
$ws1.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202409/3NbN422C1726052875488.jpeg"

$ws1.Range("B3").Value = "'2024-10-26"
$ws1.Range("C3").Value = "南宁·熊喵M动漫嘉年华·万圣派对"
$ws1.Range("D3").Value = "亭洪路45号 百益上河城"
$ws1.Range("E3").Value = "2024.10.26 11:00-10.27 21:00"
$ws1.Range("F3").Value = 120
$ws1.Range("G3").Value = 60
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=91894"
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202409/hBNwkgri1727595368161.jpeg"

$ws1.Range("B4").Value = "'2024-11-02"
$ws1.Range("C4").Value = "南宁·万圣漫控嘉年华10"
$ws1.Range("D4").Value = "亭洪路45号 百益上河城"
$ws1.Range("E4").Value = "2024.11.02 11:00-11.03 22:00"
$ws1.Range("F4").Value = 658
$ws1.Range("G4").Value = 50
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=87820"
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202409/mDTW4lHL1727235917704.jpeg"

$ws1.Range("B5").Value = "'2024-11-02"
$ws1.Range("C5").Value = "南宁·梦中礼Lolita茶会"
$ws1.Range("D5").Value = "吉兴西路盛天汇一、三、四层 云庭汇·安吉宴会厅"
$ws1.Range("E5").Value = "2024.11.02 13:00-11.02 17:00"
$ws1.Range("F5").Value = 54
$ws1.Range("G5").Value = 138
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=92826"
$ws1.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202409/09AXaAJA1726816540668.jpeg"

# Drop the now-duplicated last row; dimension A1:I6 -> A1:I5
$ws1.Rows.Item(6).Delete()

$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("B2").Value = "'2024-10-19"
$ws4.Range("C2").Value = "南宁·10.19剑网3同人only——寒光烈火·阵营PK战"
$ws4.Range("D2").Value = "大学东路158号 维也纳酒店动物园店"
$ws4.Range("E2").Value = "2024.10.19 10:00-10.19 17:30"
$ws4.Range("F2").Value = 7
$ws4.Range("G2").Value = 78
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=92730"
$ws4.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202409/3NbN422C1726052875488.jpeg"

$ws4.Range("B3").Value = "'2024-10-19"
$ws4.Range("C3").Value = "南宁·井草圣二 2024《落叶轻扬》指弹吉他音乐会"
$ws4.Range("D3").Value = "亭洪路45号 上河城艺术中心"
$ws4.Range("E3").Value = "2024.10.19 19:30-10.19 21:00"
$ws4.Range("F3").Value = 8
$ws4.Range("G3").Value = 260
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=91345"
$ws4.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202408/7rcuyrqP1724741707788.jpeg"

$ws4.Range("B4").Value = "'2024-10-26"
$ws4.Range("C4").Value = "南宁·熊喵M动漫嘉年华·万圣派对"
$ws4.Range("D4").Value = "亭洪路45号 百益上河城"
$ws4.Range("E4").Value = "2024.10.26 11:00-10.27 21:00"
$ws4.Range("F4").Value = 120
$ws4.Range("G4").Value = 60
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=91894"
$ws4.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202409/hBNwkgri1727595368161.jpeg"

$ws4.Range("B5").Value = "'2024-11-02"
$ws4.Range("C5").Value = "南宁·万圣漫控嘉年华10"
$ws4.Range("D5").Value = "亭洪路45号 百益上河城"
$ws4.Range("E5").Value = "2024.11.02 11:00-11.03 22:00"
$ws4.Range("F5").Value = 658
$ws4.Range("G5").Value = 50
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=87820"
$ws4.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202409/mDTW4lHL1727235917704.jpeg"

$ws4.Range("B6").Value = "'2024-11-02"
$ws4.Range("C6").Value = "南宁·梦中礼Lolita茶会"
$ws4.Range("D6").Value = "吉兴西路盛天汇一、三、四层 云庭汇·安吉宴会厅"
$ws4.Range("E6").Value = "2024.11.02 13:00-11.02 17:00"
$ws4.Range("F6").Value = 54
$ws4.Range("G6").Value = 138
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=92826"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202409/09AXaAJA1726816540668.jpeg"

# Drop the now-duplicated last row; dimension A1:I7 -> A1:I6
$ws4.Rows.Item(7).Delete()
